# Update the loading-times sheet with new penalty/assignment results and
# append a new row (the data now spans A1:D10 instead of A1:D9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for rows 2-10, columns A (TruckID), B (AssignedDockPosition),
# C (start_loading_time), D (end_loading_time).
$data = @(
    @(1, 1, 5, 5),
    @(2, 1, 10, 10),
    @(3, 1, 15, 15),
    @(4, 1, 20, 20),
    @(8, 1, 25, 25),
    @(6, 2, 5, 6),
    @(5, 3, 5, 5),
    @(8, 3, 10, 10),
    @(7, 4, 5, 5)
)

$row = 2
foreach ($values in $data) {
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $row++
}
